$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2307.5715
$ws.Range("I19").Value = 2207.6667
$ws.Range("J19").Value = 2382.5
$ws.Range("K19").Value = 2207.6667
$ws.Range("L19").Value = 2382.5
$ws.Range("M19").Value = -2032.6667
$ws.Range("N19").Value = -2732.5
$ws.Range("H20").Value = 574
$ws.Range("I20").Value = 574
$ws.Range("K20").Value = 574
$ws.Range("M20").Value = -344
$ws.Range("H35").Value = 574
$ws.Range("I35").Value = 574
$ws.Range("K35").Value = 574
$ws.Range("M35").Value = -195
$ws.Range("H40").Value = 1984
$ws.Range("I40").Value = 1984
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1984
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1809
$ws.Range("N40").ClearContents()
$ws.Range("H42").Value = 6024.2
$ws.Range("J42").Value = 7510.25
$ws.Range("L42").Value = 22530.75
$ws.Range("N42").Value = -22990.75
$ws.Range("H137").Value = 2620.5715
$ws.Range("I137").Value = 968
$ws.Range("K137").Value = 2904
$ws.Range("M137").Value = -354

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 158.25
$ws.Range("I5").Value = 158.25
$ws.Range("K5").Value = 158.25
$ws.Range("M5").Value = -46.25
$ws.Range("H14").Value = 25350
$ws.Range("I14").Value = 50000
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 50000
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = -49825
$ws.Range("N14").Value = -1050
$ws.Range("H16").Value = 17442.5
$ws.Range("I16").Value = 20421.2
$ws.Range("J16").Value = 2549
$ws.Range("K16").Value = 20421.2
$ws.Range("L16").Value = 2549
$ws.Range("M16").Value = -20134.2
$ws.Range("N16").Value = -3123
$ws.Range("H53").Value = 7521.5
$ws.Range("J53").Value = 43
$ws.Range("L53").Value = 43
$ws.Range("N53").Value = -1407
$ws.Range("H102").Value = 1644.25
$ws.Range("I102").Value = 1393.4286
$ws.Range("K102").Value = 1393.4286
$ws.Range("M102").Value = 228.5714
$ws.Range("H110").Value = 2438.25
$ws.Range("I110").Value = 2503
$ws.Range("J110").Value = 2330.3333
$ws.Range("K110").Value = 2503
$ws.Range("L110").Value = 2330.3333
$ws.Range("M110").Value = -458
$ws.Range("N110").Value = -6420.3333

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 158.25
$ws.Range("I4").Value = 158.25
$ws.Range("K4").Value = 158.25
$ws.Range("M4").Value = -43.25
$ws.Range("H22").Value = 105.5
$ws.Range("I22").Value = 111
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 111
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 62
$ws.Range("N22").Value = -446
$ws.Range("H54").Value = 22027.666
$ws.Range("I54").Value = 22027.666
$ws.Range("K54").Value = 22027.666
$ws.Range("M54").Value = -21543.666
$ws.Range("H80").Value = 700.5
$ws.Range("I80").Value = 694.5
$ws.Range("J80").Value = 702.5
$ws.Range("K80").Value = 694.5
$ws.Range("L80").Value = 702.5
$ws.Range("M80").Value = 303.5
$ws.Range("N80").Value = -2698.5
$ws.Range("H83").Value = 700.5
$ws.Range("I83").Value = 694.5
$ws.Range("J83").Value = 702.5
$ws.Range("K83").Value = 3472.5
$ws.Range("L83").Value = 3512.5
$ws.Range("M83").Value = 1519.5
$ws.Range("N83").Value = -13496.5
$ws.Range("H105").Value = 2626.182
$ws.Range("I105").Value = 1598.6666
$ws.Range("K105").Value = 1598.6666
$ws.Range("M105").Value = 148.3334

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 141.28572
$ws.Range("I7").Value = 137.8
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 137.8
$ws.Range("L7").Value = 150
$ws.Range("M7").Value = -24.80000000000001
$ws.Range("N7").Value = -376
$ws.Range("H16").Value = 1799.8334
$ws.Range("I16").Value = 1799.8334
$ws.Range("K16").Value = 1799.8334
$ws.Range("M16").Value = -1512.8334
$ws.Range("H31").Value = 1925.0189
$ws.Range("I31").Value = 1462.4103
$ws.Range("K31").Value = 1462.4103
$ws.Range("M31").Value = -1167.4103
$ws.Range("H34").Value = 1925.0189
$ws.Range("I34").Value = 1462.4103
$ws.Range("K34").Value = 1462.4103
$ws.Range("M34").Value = -1260.4103
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H86").Value = 8617.462
$ws.Range("I86").Value = 7713.9
$ws.Range("K86").Value = 7713.9
$ws.Range("M86").Value = -6590.9
$ws.Range("H89").Value = 8617.462
$ws.Range("I89").Value = 7713.9
$ws.Range("K89").Value = 38569.5
$ws.Range("M89").Value = -32953.5
$ws.Range("H99").Value = 6090
$ws.Range("J99").Value = 5599.5
$ws.Range("L99").Value = 5599.5
$ws.Range("N99").Value = -8595.5
$ws.Range("H113").Value = 1799.8334
$ws.Range("I113").Value = 1799.8334
$ws.Range("K113").Value = 1799.8334
$ws.Range("M113").Value = 370.1666
$ws.Range("H122").Value = 1040.4
$ws.Range("I122").Value = 925.5
$ws.Range("K122").Value = 2776.5
$ws.Range("M122").Value = -326.5
$ws.Range("H126").Value = 6090
$ws.Range("J126").Value = 5599.5
$ws.Range("L126").Value = 16798.5
$ws.Range("N126").Value = -21738.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 378
$ws.Range("I2").Value = 533
$ws.Range("J2").Value = 331.5
$ws.Range("K2").Value = 3198
$ws.Range("L2").Value = 1989
$ws.Range("M2").Value = -3085
$ws.Range("N2").Value = -2215
$ws.Range("H3").Value = 4750
$ws.Range("J3").Value = 4500
$ws.Range("L3").Value = 13500
$ws.Range("N3").Value = -13724
$ws.Range("H10").Value = 11.166667
$ws.Range("I10").Value = 11.166667
$ws.Range("K10").Value = 33.500001
$ws.Range("M10").Value = 105.499999
$ws.Range("H11").Value = 200645.1
$ws.Range("I11").Value = 400604.2
$ws.Range("J11").Value = 686
$ws.Range("K11").Value = 1201812.6
$ws.Range("L11").Value = 2058
$ws.Range("M11").Value = -1201672.6
$ws.Range("N11").Value = -2338
$ws.Range("H13").Value = 186
$ws.Range("I13").Value = 466.33334
$ws.Range("J13").Value = 17.8
$ws.Range("K13").Value = 1399.00002
$ws.Range("L13").Value = 53.40000000000001
$ws.Range("M13").Value = -1231.00002
$ws.Range("N13").Value = -389.4
$ws.Range("H38").Value = 78.111115
$ws.Range("J38").Value = 135
$ws.Range("L38").Value = 405
$ws.Range("N38").Value = -1099
$ws.Range("H103").Value = 168.42857
$ws.Range("I103").Value = 188.16667
$ws.Range("J103").Value = 50
$ws.Range("K103").Value = 564.50001
$ws.Range("L103").Value = 150
$ws.Range("M103").Value = 314.49999
$ws.Range("N103").Value = -1908
$ws.Range("H117").Value = 700.3333
$ws.Range("J117").Value = 666
$ws.Range("L117").Value = 1998
$ws.Range("N117").Value = -8882
$ws.Range("H121").Value = 875.3333
$ws.Range("I121").Value = 781.1667
$ws.Range("J121").Value = 922.4167
$ws.Range("K121").Value = 2343.5001
$ws.Range("L121").Value = 2767.2501
$ws.Range("M121").Value = -1033.5001
$ws.Range("N121").Value = -5387.2501
$ws.Range("H131").Value = 1069
$ws.Range("J131").Value = 1188.3334
$ws.Range("L131").Value = 3565.0002
$ws.Range("N131").Value = -13645.0002
$ws.Range("H139").Value = 3509.0908
$ws.Range("I139").Value = 900
$ws.Range("K139").Value = 2700
$ws.Range("M139").Value = 2440

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2799.875
$ws.Range("I102").Value = 2514.1428
$ws.Range("K102").Value = 2514.1428
$ws.Range("M102").Value = -892.1428000000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1184.4
$ws.Range("I10").Value = 1443
$ws.Range("K10").Value = 1443
$ws.Range("M10").Value = -1303
$ws.Range("H22").Value = 859.4
$ws.Range("I22").Value = 895.0833
$ws.Range("J22").Value = 716.6667
$ws.Range("K22").Value = 895.0833
$ws.Range("L22").Value = 716.6667
$ws.Range("M22").Value = -600.0833
$ws.Range("N22").Value = -1306.6667
$ws.Range("H27").Value = 859.4
$ws.Range("I27").Value = 895.0833
$ws.Range("J27").Value = 716.6667
$ws.Range("K27").Value = 895.0833
$ws.Range("L27").Value = 716.6667
$ws.Range("M27").Value = -788.0833
$ws.Range("N27").Value = -930.6667
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H136").Value = 2922
$ws.Range("I136").Value = 2922
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8766
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6216
$ws.Range("N136").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 17026.834
$ws.Range("I23").Value = 20430
$ws.Range("J23").Value = 11
$ws.Range("K23").Value = 20430
$ws.Range("L23").Value = 11
$ws.Range("M23").Value = -20201
$ws.Range("N23").Value = -469
$ws.Range("H54").Value = 22800
$ws.Range("I54").Value = 12000
$ws.Range("K54").Value = 12000
$ws.Range("M54").Value = -11480
